$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Re-bucket the "Usertype" column (E) for the existing rows (13-26) so
#    that, once the 25 new rows below are appended, each Usertype occupies a
#    contiguous block: Student (2-24), Teacher (25-42), Administrator (43-51).
# ---------------------------------------------------------------------------
$studentRows = 13..24
foreach ($r in $studentRows) {
    $ws.Range("E$r").Value = "Student"
}

$teacherRows = 25..26
foreach ($r in $teacherRows) {
    $ws.Range("E$r").Value = "Teacher"
}

# ---------------------------------------------------------------------------
# 2) Append the 25 new roster records (rows 27-51: 16 Teachers, 9 Administrators)
# ---------------------------------------------------------------------------
$newRecords = @(
    @{ Row = 27; First = "Tamara ";   Last = "Gamble";    Id = 100025; Type = "Teacher" }
    @{ Row = 28; First = "Dustin";    Last = "Phelps";    Id = 100026; Type = "Teacher" }
    @{ Row = 29; First = "Christina"; Last = "Middleton"; Id = 100027; Type = "Teacher" }
    @{ Row = 30; First = "Cali";      Last = "Pearson";   Id = 100028; Type = "Teacher" }
    @{ Row = 31; First = "Russell";   Last = "Rowe";      Id = 100029; Type = "Teacher" }
    @{ Row = 32; First = "Kyle ";     Last = "Ramos";     Id = 100030; Type = "Teacher" }
    @{ Row = 33; First = "Jayleen";   Last = "McConnell"; Id = 100031; Type = "Teacher" }
    @{ Row = 34; First = "Matthias";  Last = "Bruce";     Id = 100032; Type = "Teacher" }
    @{ Row = 35; First = "Marley";    Last = "Mooney";    Id = 100033; Type = "Teacher" }
    @{ Row = 36; First = "Holly";     Last = "Nixon";     Id = 100034; Type = "Teacher" }
    @{ Row = 37; First = "Natalie";   Last = "Woods";     Id = 100035; Type = "Teacher" }
    @{ Row = 38; First = "Jessie";    Last = "Cain";      Id = 100036; Type = "Teacher" }
    @{ Row = 39; First = "Kristen";   Last = "Duke";      Id = 100037; Type = "Teacher" }
    @{ Row = 40; First = "Melanie ";  Last = "Soto";      Id = 100038; Type = "Teacher" }
    @{ Row = 41; First = "Rylee";     Last = "Goodman";   Id = 100039; Type = "Teacher" }
    @{ Row = 42; First = "Savannah";  Last = "Velez";     Id = 100040; Type = "Teacher" }
    @{ Row = 43; First = "Damion ";   Last = "Cooley";    Id = 100041; Type = "Administrator" }
    @{ Row = 44; First = "Kellen";    Last = "Kirby";     Id = 100042; Type = "Administrator" }
    @{ Row = 45; First = "Drake";     Last = "Randolph";  Id = 100043; Type = "Administrator" }
    @{ Row = 46; First = "Owen";      Last = "Duran";     Id = 100044; Type = "Administrator" }
    @{ Row = 47; First = "Patrick";   Last = "Gould";     Id = 100045; Type = "Administrator" }
    @{ Row = 48; First = "Kyra";      Last = "Fox";       Id = 100046; Type = "Administrator" }
    @{ Row = 49; First = "Ethan";     Last = "Weber";     Id = 100047; Type = "Administrator" }
    @{ Row = 50; First = "Marie";     Last = "Boyd";      Id = 100048; Type = "Administrator" }
    @{ Row = 51; First = "Ashlie";    Last = "Kent";      Id = 100049; Type = "Administrator" }
)

foreach ($rec in $newRecords) {
    $r = $rec.Row
    $ws.Range("A$r").Value = $rec.First
    $ws.Range("B$r").Value = $rec.Last
    $ws.Range("C$r").Formula = "=CONCATENATE(A$r,B$r)"
    $ws.Range("D$r").Value = $rec.Id
    $ws.Range("E$r").Value = $rec.Type
}

# ---------------------------------------------------------------------------
# 3) Grow the table / autofilter so it covers the full A1:E51 range
# ---------------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E51"))

# ---------------------------------------------------------------------------
# 4) Leave the selection where the author ended up after typing the data
# ---------------------------------------------------------------------------
$ws.Range("E41").Select()
